$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item("Tabla 1")
$tbl = $shp.Table

# Insert a new "Febrero" row before the existing "Contratos Clientes" row (row 3)
$febRow = $tbl.Rows.Add(3)
$tbl.Rows.Item(3).Height = 33.56094488188976
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text = "Febrero"
$tbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "27-02-15"
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Text = "27-02-15"

# Update the former "Contratos Clientes" row (now row 4) to "Marzo"
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Marzo"
$tbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "31-03-15"
$tbl.Cell(4, 3).Shape.TextFrame.TextRange.Text = "31-03-15"

# Append a new "Abril" row at the end (row 5)
$abrRow = $tbl.Rows.Add(5)
$tbl.Cell(5, 1).Shape.TextFrame.TextRange.Text = "Abril"
$tbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "30-04-15"
$tbl.Cell(5, 3).Shape.TextFrame.TextRange.Text = "30-04-15"
